$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3197272419929504
$ws.Range("B1").Value = 0.2511351108551025
$ws.Range("C1").Value = 0.2173262983560562
$ws.Range("D1").Value = 0.2325874418020248
$ws.Range("E1").Value = 0.279276967048645
